# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions": columns D (Price) and
# E (Volume(1h)) are refreshed for rows 2-51. Values are plain text cells
# (no numeric NumberFormat in the sheet), so Price updates are entered with
# a leading apostrophe to force text interpretation and avoid Excel's
# automatic number conversion (e.g. "299.11" -> 299.11 numeric), keeping
# the cell type consistent with the original inline string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.841.12"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "'2.265.42"
$ws.Range("E3").Value = "  -3.76%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'299.11"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").Value = "'99.16"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.568"
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.504"
$ws.Range("E9").Value = "  -6.56%  "
$ws.Range("D10").Value = "'34.94"
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("D11").Value = "'0.0794"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").Value = "'7.01"
$ws.Range("E12").Value = "  -6.26%  "
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "'2.609.72"
$ws.Range("E14").Value = "  -3.68%  "
$ws.Range("D15").Value = "'2.267.19"
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("D16").Value = "'13.60"
$ws.Range("E16").Value = "  -5.17%  "
$ws.Range("D17").Value = "'46.814.46"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'0.795"
$ws.Range("E18").Value = "  -4.89%  "
$ws.Range("D19").Value = "'0.0₃0975"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").Value = "'12.39"
$ws.Range("E20").Value = "  -10.03%  "
$ws.Range("D21").Value = "'5.80"
$ws.Range("E21").Value = "  -6.69%  "
$ws.Range("D22").Value = "'65.62"
$ws.Range("E22").Value = "  -2.03%  "
$ws.Range("D23").Value = "'245.77"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'2.77"
$ws.Range("E24").Value = "  -7.43%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "'1.85"
$ws.Range("E26").Value = "  -7.55%  "
$ws.Range("D27").Value = "'41.25"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = "  -3.77%  "
$ws.Range("D29").Value = "'9.52"
$ws.Range("E29").Value = "  -4.30%  "
$ws.Range("D30").Value = "'19.97"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("E31").Value = "  +7.69%  "
$ws.Range("D32").Value = "'3.32"
$ws.Range("E32").Value = "  +5.08%  "
$ws.Range("D33").Value = "'145.37"
$ws.Range("E33").Value = "  -4.57%  "
$ws.Range("D34").Value = "'5.30"
$ws.Range("E34").Value = "  -8.58%  "
$ws.Range("D35").Value = "'0.0765"
$ws.Range("E35").Value = "  -6.73%  "
$ws.Range("D36").Value = "'0.111"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("D38").Value = "'15.45"
$ws.Range("E38").Value = "  +9.71%  "
$ws.Range("D39").Value = "'1.66"
$ws.Range("E39").Value = "  -10.52%  "
$ws.Range("D40").Value = "'3.81"
$ws.Range("E40").Value = "  -7.20%  "
$ws.Range("D41").Value = "'0.0295"
$ws.Range("E41").Value = "  -7.42%  "
$ws.Range("D42").Value = "'3.06"
$ws.Range("E42").Value = "  -11.81%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'93.06"
$ws.Range("E44").Value = "  +14.50%  "
$ws.Range("D45").Value = "'1.782.91"
$ws.Range("E45").Value = "  -4.20%  "
$ws.Range("D46").Value = "'1.88"
$ws.Range("E46").Value = "  -4.44%  "
$ws.Range("D47").Value = "'70.69"
$ws.Range("E47").Value = "  -4.16%  "
$ws.Range("D48").Value = "'0.183"
$ws.Range("E48").Value = "  -8.21%  "
$ws.Range("D49").Value = "'4.77"
$ws.Range("E49").Value = "  -3.69%  "
$ws.Range("D50").Value = "'94.16"
$ws.Range("E50").Value = "  -5.58%  "
$ws.Range("D51").Value = "'7.84"
$ws.Range("E51").Value = "  -2.52%  "
